# Auto-generated edit script applying cached-price refresh values
# to specific Leve rows across multiple sheets (ALC, ARM, BSM, CRP, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H75").Value = 41000
$ws.Range("J75").Value = 41000
$ws.Range("L75").Value = 41000
$ws.Range("N75").Value = -42872
$ws.Range("H78").Value = 41000
$ws.Range("J78").Value = 41000
$ws.Range("L78").Value = 123000
$ws.Range("N78").Value = -132360
$ws.Range("H112").Value = 1413.8096
$ws.Range("I112").Value = 1100
$ws.Range("K112").Value = 3300
$ws.Range("M112").Value = -2192
$ws.Range("H114").Value = 33950
$ws.Range("J114").Value = 33950
$ws.Range("L114").Value = 33950
$ws.Range("N114").Value = -42628
$ws.Range("H128").Value = 42500
$ws.Range("J128").Value = 42500
$ws.Range("L128").Value = 42500
$ws.Range("N128").Value = -52460
$ws.Range("H130").Value = 36686.668
$ws.Range("J130").Value = 36686.668
$ws.Range("L130").Value = 36686.668
$ws.Range("N130").Value = -46726.668
$ws.Range("H133").Value = 74445
$ws.Range("J133").Value = 74445
$ws.Range("L133").Value = 74445
$ws.Range("N133").Value = -84565
$ws.Range("H137").Value = 1846.8937
$ws.Range("I137").Value = 951.86206
$ws.Range("J137").Value = 3288.889
$ws.Range("K137").Value = 2855.58618
$ws.Range("L137").Value = 9866.667000000001
$ws.Range("M137").Value = -305.5861800000002
$ws.Range("N137").Value = -14966.667

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4771.625
$ws.Range("I74").Value = 4625.0713
$ws.Range("J74").Value = 4976.8
$ws.Range("K74").Value = 4625.0713
$ws.Range("L74").Value = 4976.8
$ws.Range("M74").Value = -3751.0713
$ws.Range("N74").Value = -6724.8
$ws.Range("H77").Value = 4771.625
$ws.Range("I77").Value = 4625.0713
$ws.Range("J77").Value = 4976.8
$ws.Range("K77").Value = 23125.3565
$ws.Range("L77").Value = 24884
$ws.Range("M77").Value = -18757.3565
$ws.Range("N77").Value = -33620
$ws.Range("H103").Value = 44888
$ws.Range("J103").Value = 44888
$ws.Range("L103").Value = 44888
$ws.Range("N103").Value = -47232
$ws.Range("H109").Value = 32500
$ws.Range("J109").Value = 32500
$ws.Range("L109").Value = 32500
$ws.Range("N109").Value = -35274
$ws.Range("H127").Value = 28000
$ws.Range("J127").Value = 28000
$ws.Range("L127").Value = 28000
$ws.Range("N127").Value = -37920
$ws.Range("H129").Value = 44374.75
$ws.Range("J129").Value = 44374.75
$ws.Range("L129").Value = 44374.75
$ws.Range("N129").Value = -54374.75
$ws.Range("H130").Value = 31803.625
$ws.Range("J130").Value = 31803.625
$ws.Range("L130").Value = 31803.625
$ws.Range("N130").Value = -41843.625
$ws.Range("H131").Value = 55484.25
$ws.Range("J131").Value = 55484.25
$ws.Range("L131").Value = 55484.25
$ws.Range("N131").Value = -65564.25

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H109").Value = 30000
$ws.Range("J109").Value = 30000
$ws.Range("L109").Value = 30000
$ws.Range("N109").Value = -32774
$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H125").Value = 44390
$ws.Range("J125").Value = 44390
$ws.Range("L125").Value = 44390
$ws.Range("N125").Value = -54230
$ws.Range("H126").Value = 28000
$ws.Range("J126").Value = 28000
$ws.Range("L126").Value = 28000
$ws.Range("N126").Value = -37880
$ws.Range("H130").Value = 79980
$ws.Range("J130").Value = 79980
$ws.Range("L130").Value = 79980
$ws.Range("N130").Value = -90020
$ws.Range("H135").Value = 69339.57000000001
$ws.Range("J135").Value = 69339.57000000001
$ws.Range("L135").Value = 69339.57000000001
$ws.Range("N135").Value = -79479.57000000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()
$ws.Range("H31").Value = 3736.2036
$ws.Range("I31").Value = 1804.5
$ws.Range("J31").Value = 6150.8335
$ws.Range("K31").Value = 1804.5
$ws.Range("L31").Value = 6150.8335
$ws.Range("M31").Value = -1509.5
$ws.Range("N31").Value = -6740.8335
$ws.Range("H34").Value = 3736.2036
$ws.Range("I34").Value = 1804.5
$ws.Range("J34").Value = 6150.8335
$ws.Range("K34").Value = 1804.5
$ws.Range("L34").Value = 6150.8335
$ws.Range("M34").Value = -1602.5
$ws.Range("N34").Value = -6554.8335
$ws.Range("H100").Value = 58403.332
$ws.Range("J100").Value = 58403.332
$ws.Range("L100").Value = 58403.332
$ws.Range("N100").Value = -60567.332
$ws.Range("H124").Value = 52663
$ws.Range("J124").Value = 52663
$ws.Range("L124").Value = 52663
$ws.Range("N124").Value = -57573
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H135").Value = 51497.65
$ws.Range("J135").Value = 51497.65
$ws.Range("L135").Value = 51497.65
$ws.Range("N135").Value = -61637.65

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 900
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 900
$ws.Range("M36").ClearContents()
$ws.Range("N36").Value = -1870
$ws.Range("H93").Value = 11740.333
$ws.Range("J93").Value = 11740.333
$ws.Range("L93").Value = 11740.333
$ws.Range("N93").Value = -15484.333
$ws.Range("H124").Value = 27997.273
$ws.Range("J124").Value = 27997.273
$ws.Range("L124").Value = 27997.273
$ws.Range("N124").Value = -37817.273
$ws.Range("H128").Value = 45648.1
$ws.Range("J128").Value = 45648.1
$ws.Range("L128").Value = 45648.1
$ws.Range("N128").Value = -55608.1

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 262250
$ws.Range("I74").Value = 9500
$ws.Range("J74").Value = 515000
$ws.Range("K74").Value = 9500
$ws.Range("L74").Value = 515000
$ws.Range("M74").Value = -8502
$ws.Range("N74").Value = -516996
$ws.Range("H77").Value = 262250
$ws.Range("I77").Value = 9500
$ws.Range("J77").Value = 515000
$ws.Range("K77").Value = 28500
$ws.Range("L77").Value = 1545000
$ws.Range("M77").Value = -23508
$ws.Range("N77").Value = -1554984
$ws.Range("H80").Value = 27804
$ws.Range("J80").Value = 27804
$ws.Range("L80").Value = 27804
$ws.Range("N80").Value = -30050
$ws.Range("H83").Value = 27804
$ws.Range("J83").Value = 27804
$ws.Range("L83").Value = 83412
$ws.Range("N83").Value = -94644
$ws.Range("H92").Value = 36289
$ws.Range("J92").Value = 36289
$ws.Range("L92").Value = 36289
$ws.Range("N92").Value = -41281
$ws.Range("H108").Value = 28684.25
$ws.Range("J108").Value = 28684.25
$ws.Range("L108").Value = 28684.25
$ws.Range("N108").Value = -36364.25
$ws.Range("H127").Value = 40000
$ws.Range("J127").Value = 40000
$ws.Range("L127").Value = 40000
$ws.Range("N127").Value = -49920
$ws.Range("H131").Value = 49980
$ws.Range("J131").Value = 49980
$ws.Range("L131").Value = 49980
$ws.Range("N131").Value = -60060
$ws.Range("H139").Value = 76635.71000000001
$ws.Range("J139").Value = 76635.71000000001
$ws.Range("L139").Value = 76635.71000000001
$ws.Range("N139").Value = -86915.71000000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 26197.25
$ws.Range("J93").Value = 26197.25
$ws.Range("L93").Value = 26197.25
$ws.Range("N93").Value = -31189.25
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H109").Value = 20968
$ws.Range("J109").Value = 20968
$ws.Range("L109").Value = 20968
$ws.Range("N109").Value = -23742
$ws.Range("H127").Value = 42607.25
$ws.Range("J127").Value = 42607.25
$ws.Range("L127").Value = 42607.25
$ws.Range("N127").Value = -52527.25
$ws.Range("H135").Value = 70522
$ws.Range("J135").Value = 70522
$ws.Range("L135").Value = 70522
$ws.Range("N135").Value = -80662
